$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Precondition text (used in TC1, TC2, TC3 "Precondition" rows B8, B19, B29)
$oldPrecondition = "Administrador esta autenticado no sistema; e, tem permissao para alterar Gerente de Desempenho"
$newPrecondition = "Administrador esta autenticado no sistema e tem permissao para alterar Gerente de Desempenho"

# 2) Expected result for step 1 (D10, D21, D31)
$oldStep1Result = "SYSTEM exibe a listagem do Catalogo (Perfis) de Competencias cadastradas com a opcao 'Alterar Gerente' dentre as varias listadas"
$newStep1Result = "SYSTEM exibe a listagem dos Perfis de Competencias cadastrados com a opcao 'Alterar Gerente' dentre as varias exibidas"

# 3) Step 3 description (B12, B23, B33)
$oldStep3Desc = "Administrador preenche o campo 'Login do Novo Gerente de Desempenho' do novo Gerente de Desempenho para o Perfil de Competencias "
$newStep3Desc = "Administrador preenche o campo 'Login do Novo Gerente de Desempenho' para o Perfil de Competencias"

# 4) TC2 step 4 expected result (D24)
$oldTc2Step4Result = "SYSTEM apresenta o Catalogo (Perfis) de Competencias cadastradas sem nenhuma alteracao"
$newTc2Step4Result = "SYSTEM apresenta o Catalogo (Perfis) de Competencias sem nenhuma alteracao"

function Replace-IfMatches {
    param(
        [string]$CellAddress,
        [string]$OldText,
        [string]$NewText
    )
    $cell = $ws.Range($CellAddress)
    if ($cell.Value2 -eq $OldText) {
        $cell.Value = $NewText
    }
}

# Precondition appears identically for TC1 (row 8), TC2 (row 19), TC3 (row 29)
foreach ($row in 8, 19, 29) {
    Replace-IfMatches "B$row" $oldPrecondition $newPrecondition
}

# Step 1 expected result appears identically for TC1 (row 10), TC2 (row 21), TC3 (row 31)
foreach ($row in 10, 21, 31) {
    Replace-IfMatches "D$row" $oldStep1Result $newStep1Result
}

# Step 3 description appears identically for TC1 (row 12), TC2 (row 23), TC3 (row 33)
foreach ($row in 12, 23, 33) {
    Replace-IfMatches "B$row" $oldStep3Desc $newStep3Desc
}

# TC2's step 4 expected result (row 24) only
Replace-IfMatches "D24" $oldTc2Step4Result $newTc2Step4Result
